$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each inner array: Row, Col, OldText, NewText
# NOTE: order matters. The cell (13,2) old value "556÷8=" is the same text
# that cell (9,5) is about to become after its own replacement. If (9,5) is
# processed first, a duplicate "556÷8=" would (temporarily) exist earlier in
# the document than (13,2)'s position, and a Find scoped to (13,2) could
# still latch onto the earlier, already-edited occurrence. To avoid that
# collision we replace (13,2) before (9,5).
$changes = @(
    @(1, 1, "761÷4=", "752÷8="),
    @(1, 2, "976÷9=", "201÷4="),
    @(1, 3, "939÷5=", "844÷8="),
    @(1, 4, "197÷3=", "338÷8="),
    @(1, 5, "409÷4=", "401÷2="),

    @(5, 1, "740÷5=", "444÷8="),
    @(5, 2, "827÷8=", "880÷4="),
    @(5, 3, "183÷9=", "135÷5="),
    @(5, 4, "308÷5=", "675÷9="),
    @(5, 5, "875÷4=", "680÷6="),

    @(9, 1, "708÷6=", "292÷6="),
    @(9, 2, "394÷5=", "468÷4="),
    @(9, 3, "280÷4=", "325÷3="),
    @(9, 4, "306÷9=", "188÷4="),

    @(13, 1, "833÷3=", "853÷5="),
    @(13, 2, "556÷8=", "847÷4="),
    @(13, 3, "315÷2=", "546÷8="),
    @(13, 4, "446÷5=", "113÷8="),
    @(13, 5, "214÷7=", "587÷2="),

    @(9, 5, "959÷3=", "556÷8="),

    @(17, 1, "885÷6=", "386÷6="),
    @(17, 2, "926÷6=", "324÷9="),
    @(17, 3, "129÷2=", "270÷4="),
    @(17, 4, "249÷7=", "107÷7="),
    @(17, 5, "355÷9=", "487÷2=")
)

foreach ($change in $changes) {
    $row = $change[0]
    $col = $change[1]
    $oldText = $change[2]
    $newText = $change[3]

    # Re-fetch the table/cell each time in case prior edits made earlier
    # handles stale.
    $t = $d.Tables.Item(1)
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 0, $false, $newText, 2)
}

Write-Host "Done"
